{"js": "// The \"Communication\" rubric header paragraph ends with a run of a single\n// space followed by a separate run of 62 spaces (63 trailing spaces total,\n// split across two runs). The edit collapses that trailing whitespace down\n// to 60 spaces living in a single run.\nconst body = context.document.body;\n\nconst headings = body.search(\"Communication\", { matchCase: true, matchWholeWord: false });\nheadings.load(\"items\");\nawait context.sync();\n\nif (headings.items.length === 0) {\n  throw new Error(\"Could not find the 'Communication' heading run\");\n}\n\nconst heading = headings.items[0];\nconst paragraph = heading.paragraphs.getFirst();\nconst paragraphEnd = paragraph.getRange(\"End\");\n\n// Range covering every character after \"Communication\" through the end of\n// the paragraph - i.e. the two whitespace-only runs (\" \" + 62 spaces).\nconst trailingSpacesRange = heading.getRange(\"After\").expandTo(paragraphEnd);\ntrailingSpacesRange.load(\"text\");\nawait context.sync();\n\n// Sanity-check we grabbed only whitespace before rewriting it.\nif (/\\S/.test(trailingSpacesRange.text)) {\n  throw new Error(\"Unexpected non-whitespace content after 'Communication'\");\n}\n\n// Replace the combined whitespace (1 + 62 = 63 spaces) with 60 spaces,\n// merging both runs into a single run.\ntrailingSpacesRange.insertText(\" \".repeat(60), \"Replace\");\nawait context.sync();\n", "ps1": "# The \"Communication\" rubric header paragraph ends with a run of a single\n# space followed by a separate run of 62 spaces (63 trailing spaces total,\n# split across two runs). The edit collapses that trailing whitespace down\n# to 60 spaces. Trimming characters off the END of the whitespace (rather\n# than rewriting the whole trailing range's .Text) keeps the still-present\n# whitespace anchored in its original run, so the run that is *not* fully\n# consumed (and the \"Communication\" run before it) keep their identity\n# instead of being coalesced into a brand new run.\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"Communication\"\n$rng.Find.MatchCase = $true\n$rng.Find.Execute() | Out-Null\n\nif (-not $rng.Find.Found) {\n  throw \"Could not find the 'Communication' heading run\"\n}\n\n$para = $rng.Paragraphs(1)\n$paraEnd = $para.Range.End\n\n# Range covering every character after \"Communication\" through the end of\n# the paragraph, excluding the paragraph mark itself (Range.End - 1).\n$trailing = $d.Range($rng.End, $paraEnd - 1)\n\nif ($trailing.Text -match \"\\S\") {\n  throw \"Unexpected non-whitespace content after 'Communication'\"\n}\n\n$targetSpaces = 60\n$currentLen = $trailing.Text.Length\n$toDelete = $currentLen - $targetSpaces\n\nif ($toDelete -gt 0) {\n  # Trim the extra spaces off the tail end (immediately before the\n  # paragraph mark) so the remaining spaces stay inside their existing run.\n  $delRange = $d.Range($paraEnd - 1 - $toDelete, $paraEnd - 1)\n  $delRange.Delete()\n}\n"}
